# "nnuevas funciones y main creado"
# Adds a demo login row to the "login" sheet and a matching account row to
# the "accounts" sheet, then leaves the "accounts" sheet as the active tab.

$wb = $excel.ActiveWorkbook

$wsLogin = $wb.Worksheets.Item("login")
$wsAccounts = $wb.Worksheets.Item("accounts")

# New row 2 on "login": a sample/admin login entry.
$wsLogin.Range("A2").Value = "66admin66"
$wsLogin.Range("B2").Value = "admin"
$wsLogin.Range("C2").Value = 1

# New row 2 on "accounts": the matching account/owner record.
$wsAccounts.Range("A2").Value = "0cc447ac8b9791feaddb339c7a63256a"
$wsAccounts.Range("B2").Value = 0
$wsAccounts.Range("C2").Value = "66admin66"

# Update the remembered selection on each sheet...
$wsLogin.Range("C2").Select()
$wsAccounts.Range("C3").Select()

# ...and make "accounts" the active/visible tab when the workbook re-opens.
$wsAccounts.Activate()
